$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New gather-event drop rows (Id, ~Name, Items, Count), following the same
# layout as the existing row (23000001 / 水池 / ...).
# Data was entered column-by-column: Id, then Items, then ~Name, then Count.

$ws.Range("A5").Value = 23000002
$ws.Range("A6").Value = 23000003
$ws.Range("A7").Value = 23000004

$ws.Range("C5").Value = "22010201;35|22010202;10|22010203;5|22010204;5"
$ws.Range("C6").Value = "22010401;30|22010402;15|22010403;10"
$ws.Range("C7").Value = "22010501;30|22010502;15|22010302;8"

$ws.Range("B5").Value = "沙堆"
$ws.Range("B6").Value = "蘑菇"
$ws.Range("B7").Value = "枯木"

# Carry the Count column style from the existing row down to the new rows,
# then set the values.
$ws.Range("D4").Copy()
$ws.Range("D5:D7").PasteSpecial(-4122)
$ws.Range("D5").Value = 3
$ws.Range("D6").Value = 3
$ws.Range("D7").Value = 3

# Grow the table ("表2") to cover the newly added rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:D7"))

# Leave the selection on the newly entered names, matching the author's
# last interactive action.
$null = $ws.Range("B5:B7").Select()
